$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of incident log data (2024-05-14 shift) to append below the
# existing table, matching the source data exactly (plain text, not
# auto-converted dates/times/numbers).
$newRows = @(
  @("2024-05-14","09:30:16","No coge placa","-","-","-","-"),
  @("2024-05-14","09:30:23","AOI no detecta pieza","-","-","-","-"),
  @("2024-05-14","09:30:31","No coge placa","-","-","-","-"),
  @("2024-05-14","09:30:38","Fallo cámara visión","-","-","-","-"),
  @("2024-05-14","09:30:44","-","-","Etiquetadora","-","-"),
  @("2024-05-14","09:30:53","-","-","-","No coloca bien la pcb","-"),
  @("2024-05-14","09:30:56","-","-","-","Fallo cámara cover","-"),
  @("2024-05-14","09:34:15","-","-","Power atascado en prensa, cuesta sacar","-","-"),
  @("2024-05-14","09:51:40","Secuencia atornillador","-","-","-","-"),
  @("2024-05-14","09:56:09","No pone tornillo","-","-","-","-"),
  @("2024-05-14","09:56:59","-","-","-","Fallo visión core","-"),
  @("2024-05-14","09:57:04","-","-","-","Repeat funcional","-"),
  @("2024-05-14","09:57:08","-","-","-","Fallo cámara ferrite","-")
)

$startRow = 115
$endRow = $startRow + $newRows.Count - 1

# Write every value through a literal-text formula first (="...") so that
# date-looking / time-looking strings such as "2024-05-14" or "09:30:16"
# are not auto-parsed into date/time serials by the usual Value-assignment
# smart-parsing. Then convert the whole block to static values in-place via
# Copy + PasteSpecial (values only), which keeps the default "General"
# cell style/format intact (same as the rest of the sheet).
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $row = $startRow + $i
  $rowData = $newRows[$i]
  for ($c = 0; $c -lt $rowData.Count; $c++) {
    $colLetter = [char](65 + $c)
    $escaped = $rowData[$c] -replace '"', '""'
    $ws.Range("$colLetter$row").Formula = "=""$escaped"""
  }
}

$newRange = $ws.Range("A$($startRow):G$($endRow)")
$newRange.Copy()
$newRange.PasteSpecial(-4163) | Out-Null
